$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update description (B2), drop the old C2 value ---
$ws.Range("B2").Value = "Карандаш зеленый обновление"
$ws.Range("C2").ClearContents()

# --- Row 3: drop the old B3 value, update description (C3) ---
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "Самые клевые ручки на свете. обновление"

# --- Row 4: update both description cells ---
$ws.Range("B4").Value = "Маркер обновление"
$ws.Range("C4").Value = "Просто уникальные и самые распрекрасные маркеры. обновление"

# --- Row 5 (new): error / test row ---
$ws.Range("A5").Value = "error code"
$ws.Range("B5").Value = "тест на ошибку"
$ws.Range("C5").Value = "тест на ошибку"

# --- View state: zoom level and active selection ---
$excel.ActiveWindow.Zoom = 145
$ws.Range("A2").Select() | Out-Null
